$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change the D1 header text from "Telefone Celular" to "codigo de barras"
$ws.Range("D1").Value = "codigo de barras"

# 2. Add a new empty styled cell at D2 (bold font style, same as header style but applied to a blank cell)
$ws.Range("D2").Value = ""
$ws.Range("D2").Font.Bold = $true

# 3. Move the active selection to D2
$ws.Range("D2").Select()

# 4. Set up the page setup (paper size + orientation) which will register printer settings
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
